$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at the front. This shifts the existing "Session"
# column (A, with its custom width) into column B, and the existing
# "Time" column (B) into column C - which is exactly the layout the
# new table needs, and it preserves the original custom column width
# metadata on the (now) second column without us having to recreate it.
$ws.Columns.Item(1).Insert()

# New table: Areas | Topics | Time
$ws.Cells.Item(1, 1).Value = "Areas"
$ws.Cells.Item(1, 2).Value = "Topics"
$ws.Cells.Item(1, 3).Value = "Time"

$ws.Cells.Item(2, 1).Value = "Domain"
$ws.Cells.Item(2, 2).Value = "Basic Biology"
$ws.Cells.Item(2, 3).Value = "1 hr"

$ws.Cells.Item(3, 1).Value = "Domain"
$ws.Cells.Item(3, 2).Value = "caGrid Induction"
$ws.Cells.Item(3, 3).Value = "2 hr"

$ws.Cells.Item(4, 1).Value = "Product Demo"
$ws.Cells.Item(4, 2).Value = "Admin Demo"
$ws.Cells.Item(4, 3).Value = "1 hr"

$ws.Cells.Item(5, 1).Value = "Product Demo"
$ws.Cells.Item(5, 2).Value = "Thick Client Demo"
$ws.Cells.Item(5, 3).Value = "2 hr"

$ws.Cells.Item(6, 1).Value = "Product Demo"
$ws.Cells.Item(6, 2).Value = "Web app Demo"
$ws.Cells.Item(6, 3).Value = "1 hr"

$ws.Cells.Item(7, 1).Value = "Technical Session"
$ws.Cells.Item(7, 2).Value = "Overall architechture"
$ws.Cells.Item(7, 3).Value = "3 hr"

$ws.Cells.Item(8, 1).Value = "Technical Session"
$ws.Cells.Item(8, 2).Value = "Introduction to model"
$ws.Cells.Item(8, 3).Value = "3 hr"

$ws.Cells.Item(9, 1).Value = "General"
$ws.Cells.Item(9, 2).Value = "Doubts Clarification"
$ws.Cells.Item(9, 3).Value = "30 min"

# Bold header row
$ws.Range("A1:C1").Font.Bold = $true

# Keep print orientation explicit (portrait), as in the updated sheet
$ws.PageSetup.Orientation = 1

# Match the saved selection/active cell from the author's edit
$ws.Range("B12").Select()
